$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported between the existing rows that
# were at 80 and 81 (old row 81, date 2020-11-27). Insert a new row at 81,
# shifting the former rows 81-84 down to 82-85, then populate the new row
# with the latest "Primera" quality record (2022-01-24).
$ws.Rows.Item(81).Insert()

$ws.Cells.Item(81, 1).Value = 10
$ws.Cells.Item(81, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(81, 3).Value = "La Araucanía"
$ws.Cells.Item(81, 4).Value = 44585
$ws.Cells.Item(81, 5).Value = 9
$ws.Cells.Item(81, 6).Value = "Fruta"
$ws.Cells.Item(81, 7).Value = 100101
$ws.Cells.Item(81, 8).Value = "Berries"
$ws.Cells.Item(81, 9).Value = 100101001
$ws.Cells.Item(81, 10).Value = "Arándano (blue)"
$ws.Cells.Item(81, 11).Value = "Sin especificar"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 200
$ws.Cells.Item(81, 14).Value = 2000
$ws.Cells.Item(81, 15).Value = 2000
$ws.Cells.Item(81, 16).Value = 2000
$ws.Cells.Item(81, 17).Value = "$/kilo"
$ws.Cells.Item(81, 18).Value = "Región del Maule"
$ws.Cells.Item(81, 19).Value = 2000
$ws.Cells.Item(81, 20).Value = 1

# Give the new date cell the same date/time number format as the rest of
# column D (style index 2 in the original workbook).
$ws.Cells.Item(81, 4).NumberFormat = $ws.Cells.Item(82, 4).NumberFormat
